# Apply scheduled-runner profit recalculations to Sheets/Phoenix_Profits.xlsx
# (currentAveragePrice / LevePrice / LeveProfit columns per-sheet, per the diff)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H59").Value = 3246.8572
$ws.Range("I59").Value = 925
$ws.Range("K59").Value = 2775
$ws.Range("M59").Value = -2218
$ws.Range("H96").Value = 15352.579
$ws.Range("I96").Value = 29663.111
$ws.Range("J96").Value = 2473.1
$ws.Range("K96").Value = 88989.333
$ws.Range("L96").Value = 7419.299999999999
$ws.Range("M96").Value = -87616.333
$ws.Range("N96").Value = -10165.3
$ws.Range("H103").Value = 2836.25
$ws.Range("I103").Value = 759.6667
$ws.Range("J103").Value = 3315.4614
$ws.Range("K103").Value = 2279.0001
$ws.Range("L103").Value = 9946.3842
$ws.Range("M103").Value = -1693.0001
$ws.Range("N103").Value = -11118.3842
$ws.Range("H130").Value = 84999.664
$ws.Range("J130").Value = 84999.664
$ws.Range("L130").Value = 84999.664
$ws.Range("N130").Value = -95039.664
$ws.Range("H132").Value = 2538.7058
$ws.Range("I132").Value = 2080.5334
$ws.Range("K132").Value = 6241.600199999999
$ws.Range("M132").Value = -3711.600199999999
$ws.Range("H136").Value = 71782.86
$ws.Range("J136").Value = 71782.86
$ws.Range("L136").Value = 71782.86
$ws.Range("N136").Value = -81982.86

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 2445.5264
$ws.Range("J4").Value = 5832.6665
$ws.Range("L4").Value = 5832.6665
$ws.Range("N4").Value = -6064.6665
$ws.Range("H32").Value = 21958.262
$ws.Range("I32").Value = 25099.447
$ws.Range("J32").Value = 15247.546
$ws.Range("K32").Value = 25099.447
$ws.Range("L32").Value = 15247.546
$ws.Range("M32").Value = -24812.447
$ws.Range("N32").Value = -15821.546
$ws.Range("H45").Value = 1616.6428
$ws.Range("I45").Value = 1135.8334
$ws.Range("K45").Value = 1135.8334
$ws.Range("M45").Value = -758.8334
$ws.Range("H74").Value = 108213.94
$ws.Range("I74").Value = 83444.27
$ws.Range("K74").Value = 83444.27
$ws.Range("M74").Value = -82570.27
$ws.Range("H77").Value = 108213.94
$ws.Range("I77").Value = 83444.27
$ws.Range("K77").Value = 417221.35
$ws.Range("M77").Value = -412853.35

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 18499.148
$ws.Range("I107").Value = 19739.12
$ws.Range("K107").Value = 19739.12
$ws.Range("M107").Value = -17819.12

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2617.85
$ws.Range("I31").Value = 2196.926
$ws.Range("J31").Value = 3492.077
$ws.Range("K31").Value = 2196.926
$ws.Range("L31").Value = 3492.077
$ws.Range("M31").Value = -1901.926
$ws.Range("N31").Value = -4082.077
$ws.Range("H34").Value = 2617.85
$ws.Range("I34").Value = 2196.926
$ws.Range("J34").Value = 3492.077
$ws.Range("K34").Value = 2196.926
$ws.Range("L34").Value = 3492.077
$ws.Range("M34").Value = -1994.926
$ws.Range("N34").Value = -3896.077
$ws.Range("H122").Value = 2217.5
$ws.Range("I122").Value = 1811.2
$ws.Range("K122").Value = 5433.6
$ws.Range("M122").Value = -2983.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1417.7059
$ws.Range("I12").Value = 326.33334
$ws.Range("J12").Value = 1651.5714
$ws.Range("K12").Value = 979.0000200000001
$ws.Range("L12").Value = 4954.7142
$ws.Range("M12").Value = -806.0000200000001
$ws.Range("N12").Value = -5300.7142
$ws.Range("H23").Value = 777.4
$ws.Range("I23").Value = 961.25
$ws.Range("J23").Value = 654.8333
$ws.Range("K23").Value = 2883.75
$ws.Range("L23").Value = 1964.4999
$ws.Range("M23").Value = -2648.75
$ws.Range("N23").Value = -2434.4999
$ws.Range("H57").Value = 11499.333
$ws.Range("I57").Value = 11499.333
$ws.Range("K57").Value = 34497.999
$ws.Range("M57").Value = -33938.999
$ws.Range("H63").Value = 7750
$ws.Range("J63").Value = 7500
$ws.Range("L63").Value = 22500
$ws.Range("N63").Value = -23998
$ws.Range("H66").Value = 7750
$ws.Range("J66").Value = 7500
$ws.Range("L66").Value = 67500
$ws.Range("N66").Value = -74988
$ws.Range("H131").Value = 4924.273
$ws.Range("I131").Value = 2295
$ws.Range("J131").Value = 6426.7144
$ws.Range("K131").Value = 6885
$ws.Range("L131").Value = 19280.1432
$ws.Range("M131").Value = -1845
$ws.Range("N131").Value = -29360.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1071.3125
$ws.Range("I97").Value = 846
$ws.Range("J97").Value = 1747.25
$ws.Range("K97").Value = 846
$ws.Range("L97").Value = 1747.25
$ws.Range("M97").Value = -350
$ws.Range("N97").Value = -2739.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 49969.895
$ws.Range("I40").Value = 59409.258
$ws.Range("K40").Value = 59409.258
$ws.Range("M40").Value = -59273.258
$ws.Range("H61").Value = 3481.1365
$ws.Range("I61").Value = 2769
$ws.Range("J61").Value = 10602.5
$ws.Range("K61").Value = 2769
$ws.Range("L61").Value = 10602.5
$ws.Range("M61").Value = -2567
$ws.Range("N61").Value = -11006.5
$ws.Range("H68").Value = 4524.1875
$ws.Range("I68").Value = 3199.5
$ws.Range("J68").Value = 6732
$ws.Range("K68").Value = 3199.5
$ws.Range("L68").Value = 6732
$ws.Range("M68").Value = -2450.5
$ws.Range("N68").Value = -8230
$ws.Range("H71").Value = 4524.1875
$ws.Range("I71").Value = 3199.5
$ws.Range("J71").Value = 6732
$ws.Range("K71").Value = 15997.5
$ws.Range("L71").Value = 33660
$ws.Range("M71").Value = -12253.5
$ws.Range("N71").Value = -41148
$ws.Range("H93").Value = 3117.6
$ws.Range("I93").Value = 2922.75
$ws.Range("K93").Value = 2922.75
$ws.Range("M93").Value = -1674.75
$ws.Range("H113").Value = 3481.1365
$ws.Range("I113").Value = 2769
$ws.Range("J113").Value = 10602.5
$ws.Range("K113").Value = 2769
$ws.Range("L113").Value = 10602.5
$ws.Range("M113").Value = -599
$ws.Range("N113").Value = -14942.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 5726.5
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 5726.5
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 5726.5
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -6062.5
$ws.Range("H41").Value = 18321
$ws.Range("J41").Value = 18321
$ws.Range("L41").Value = 18321
$ws.Range("N41").Value = -19101
$ws.Range("H114").Value = 50000
$ws.Range("J114").Value = 50000
$ws.Range("L114").Value = 50000
$ws.Range("N114").Value = -58678
$ws.Range("H122").Value = 2287.2058
$ws.Range("I122").Value = 2273.0645
$ws.Range("K122").Value = 6819.193499999999
$ws.Range("M122").Value = -4369.193499999999
